$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr6 = New-Object 'object[,]' 1,16
$arr6[0,0] = [double]"2.8485153009081172E-2"
$arr6[0,1] = [double]"0.86297334996570385"
$arr6[0,2] = [double]"2.6344747351444719E-2"
$arr6[0,3] = [double]"220283961.18990496"
$arr6[0,4] = [double]"27.385372260478452"
$arr6[0,5] = [double]"62.484846405793988"
$arr6[0,6] = [double]"5.1341982211688406E-3"
$arr6[0,7] = [double]"1.0689880108764041E-5"
$arr6[0,8] = [double]"1.70599990405285E-4"
$arr6[0,9] = [double]"5.7371436298185564E-4"
$arr6[0,10] = [double]"0.99284636329150489"
$arr6[0,11] = [double]"0.9961313371698316"
$arr6[0,12] = [double]"0.96572359980665112"
$arr6[0,13] = [double]"8.4155215770682495E-2"
$arr6[0,14] = [double]"6.1886801515886339E-2"
$arr6[0,15] = [double]"0.1842108471057429"
$ws.Range("A6:P6").Value = $arr6

$arr13 = New-Object 'object[,]' 1,16
$arr13[0,0] = [double]"1.8161402303727852E-2"
$arr13[0,1] = [double]"0.29768615002318066"
$arr13[0,2] = [double]"2.6815875178212798E-2"
$arr13[0,3] = [double]"492855589.35860401"
$arr13[0,4] = [double]"33.011665073474873"
$arr13[0,5] = [double]"63.143880638923228"
$arr13[0,6] = [double]"1.1279517962478399E-3"
$arr13[0,7] = [double]"9.0448712362826405E-3"
$arr13[0,8] = [double]"0.60102695522390481"
$arr13[0,9] = [double]"2.8734711070142958E-6"
$arr13[0,10] = [double]"0.99779842035800714"
$arr13[0,11] = [double]"0.98224323259105983"
$arr13[0,12] = [double]"0.7564785955631903"
$arr13[0,13] = [double]"4.6685799185329745E-2"
$arr13[0,14] = [double]"0.13258657448946634"
$arr13[0,15] = [double]"0.49100528550356926"
$ws.Range("A13:P13").Value = $arr13

$arr15 = New-Object 'object[,]' 1,16
$arr15[0,0] = [double]"2.4964693400963234E-2"
$arr15[0,1] = [double]"0.35337765057996651"
$arr15[0,2] = [double]"4.1409677904350399E-2"
$arr15[0,3] = [double]"505562080.06205565"
$arr15[0,4] = [double]"23.311738782413769"
$arr15[0,5] = [double]"73.531399846732299"
$arr15[0,6] = [double]"3.6134982773081098E-3"
$arr15[0,7] = [double]"0.12884805991454981"
$arr15[0,8] = [double]"1.5890966837976868E-3"
$arr15[0,9] = [double]"4.1593658086928869E-3"
$arr15[0,10] = [double]"0.99724616386640075"
$arr15[0,11] = [double]"0.99443910669687763"
$arr15[0,12] = [double]"0.77735570342480309"
$arr15[0,13] = [double]"5.2213961468780483E-2"
$arr15[0,14] = [double]"7.4197603533342857E-2"
$arr15[0,15] = [double]"0.4694867981205062"
$ws.Range("A15:P15").Value = $arr15
